# Timing issue fix - keywords, updated tc1,2 in ubc01
#
# The "CasesTab" query (row 2, column B) previously returned an extra
# `Cohort` column (via coalesce(co.cohort_description, '') AS `Cohort`).
# That trailing clause (and the now-dangling comma after the
# "Response to Treatment" column) is removed here so the query once again
# ends on the `Response to Treatment` column.
#
# The row also shrinks slightly to reflect the now-shorter query text, and
# the sheet's saved view/selection is moved back up to the top of the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the CasesTab query text in B2: drop the trailing Cohort column ---
$marker = "coalesce(diag.best_response, '') AS ``Response to Treatment``"
$oldB2 = $ws.Range("B2").Value()
$idx = $oldB2.IndexOf($marker)
if ($idx -ge 0) {
    $newB2 = $oldB2.Substring(0, $idx) + $marker
    $ws.Range("B2").Value = $newB2
}

# --- Row heights: row 2 shrinks now that its query text is shorter ---
$ws.Rows(2).RowHeight = 290

# --- Restore the saved view to the top of the sheet/selection ---
$ws.Range("B2").Select()
